$d = $word.ActiveDocument

# Locate the last paragraph in the body (list item ending the document)
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
$lastRange.Collapse(0)

# Create a brand-new paragraph after it; InsertParagraphAfter clones the
# owning paragraph's pPr (pStyle "Akapitzlist" + numPr ilvl=0/numId=3).
$lastRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range
$newRange.Collapse(0)
$newText = "Podpowiedzi gdy chcemy napisać do kogoś wiadomość (>, >>), tak aby można było wybrać osobę, którą się obserwuje takie intelisence "
$newRange.InsertAfter($newText)

# Work out the character offsets of the new (still-empty-of-sym) paragraph
# so we can target the boundary between it and a helper paragraph below.
$paraStart = $newPara.Range.Start
$paraTextEnd = $paraStart + $newText.Length

# InsertXML always lands its payload as a *new sibling paragraph* right
# after the paragraph owning the target range, so give the helper
# paragraph the very same pPr as the real target paragraph: once we
# delete the paragraph mark between them below, the merged paragraph
# keeps this pPr (the pPr of the *second*/surviving paragraph wins on a
# paragraph-mark delete, exactly like typing Delete at the end of a
# Word list item).
$endOfTextRange = $d.Range($paraTextEnd, $paraTextEnd)
$symXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Akapitzlist"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:sym w:font="Wingdings" w:char="F04A"/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$endOfTextRange.InsertXML($symXml)

# Merge the helper paragraph back into the target paragraph by deleting
# the paragraph mark that currently separates them, leaving a single
# paragraph whose second run is the Wingdings smiley symbol.
$joinRange = $d.Range($paraTextEnd, $paraTextEnd + 1)
$joinRange.Delete()
